$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "30.004.04"
Set-TextValue $ws "E2" "  -0.29%  "
Set-TextValue $ws "D3" "1.868.85"
Set-TextValue $ws "E3" "  -2.83%  "
Set-TextValue $ws "E4" "  +0.12%  "
Set-TextValue $ws "D5" "317.92"
Set-TextValue $ws "E5" "  -2.43%  "
Set-TextValue $ws "E6" "  +0.06%  "
Set-TextValue $ws "D7" "0.5078"
Set-TextValue $ws "E7" "  -1.63%  "
Set-TextValue $ws "D8" "0.3907"
Set-TextValue $ws "E8" "  -2.35%  "
Set-TextValue $ws "D9" "0.08141"
Set-TextValue $ws "E9" "  -3.91%  "
Set-TextValue $ws "D10" "41.97"
Set-TextValue $ws "E10" "  -2.20%  "
Set-TextValue $ws "D11" "1.089"
Set-TextValue $ws "E11" "  -3.03%  "
Set-TextValue $ws "D12" "22.59"
Set-TextValue $ws "E12" "  +6.65%  "
Set-TextValue $ws "D13" "1.872.91"
Set-TextValue $ws "E13" "  -2.50%  "
Set-TextValue $ws "D14" "6.248"
Set-TextValue $ws "E14" "  -1.41%  "
Set-TextValue $ws "D15" "7.143"
Set-TextValue $ws "E15" "  -2.74%  "
Set-TextValue $ws "E16" "  +0.12%  "
Set-TextValue $ws "D17" "91.52"
Set-TextValue $ws "E17" "  -3.01%  "
Set-TextValue $ws "D18" "0.00001075"
Set-TextValue $ws "E18" "  -3.64%  "
Set-TextValue $ws "D19" "0.06331"
Set-TextValue $ws "E19" "  -6.37%  "
Set-TextValue $ws "D20" "17.82"
Set-TextValue $ws "E20" "  -1.08%  "
Set-TextValue $ws "E21" "  +0.02%  "
Set-TextValue $ws "D22" "29.996.42"
Set-TextValue $ws "E22" "  -0.29%  "
Set-TextValue $ws "D23" "5.777"
Set-TextValue $ws "E23" "  -4.68%  "
Set-TextValue $ws "E24" "  -1.25%  "
Set-TextValue $ws "D25" "2.201"
Set-TextValue $ws "D26" "2.088.37"
Set-TextValue $ws "E26" "  -2.48%  "
Set-TextValue $ws "D27" "160.37"
Set-TextValue $ws "E27" "  +0.24%  "
Set-TextValue $ws "D28" "20.83"
Set-TextValue $ws "E28" "  -0.76%  "
Set-TextValue $ws "D29" "2.219"
Set-TextValue $ws "E29" "  -10.02%  "
Set-TextValue $ws "D30" "126.14"
Set-TextValue $ws "E30" "  -2.33%  "
Set-TextValue $ws "E31" "  -2.74%  "
Set-TextValue $ws "D32" "1.038"
Set-TextValue $ws "E32" "  -3.84%  "
Set-TextValue $ws "D33" "5.853"
Set-TextValue $ws "E33" "  -3.64%  "
Set-TextValue $ws "D34" "3.733"
Set-TextValue $ws "E34" "  +2.08%  "
Set-TextValue $ws "D35" "0.02411"
Set-TextValue $ws "E35" "  -3.63%  "
Set-TextValue $ws "B36" "Hedera"
Set-TextValue $ws "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D36" "0.06311"
Set-TextValue $ws "E36" "  -4.40%  "
Set-TextValue $ws "B37" "InternetComputer(DFINITY)"
Set-TextValue $ws "C37" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D37" "5.149"
Set-TextValue $ws "E37" "  -1.00%  "
Set-TextValue $ws "D38" "0.2133"
Set-TextValue $ws "E38" "  -4.14%  "
Set-TextValue $ws "D39" "1.169"
Set-TextValue $ws "E39" "  -6.12%  "
Set-TextValue $ws "D40" "8.448"
Set-TextValue $ws "E40" "  -6.29%  "
Set-TextValue $ws "D41" "0.6253"
Set-TextValue $ws "E41" "  -4.38%  "
Set-TextValue $ws "D42" "1.208"
Set-TextValue $ws "E42" "  -2.62%  "
Set-TextValue $ws "D43" "11.19"
Set-TextValue $ws "E43" "  -1.80%  "
Set-TextValue $ws "D44" "0.9998"
Set-TextValue $ws "E44" "  -0.15%  "
Set-TextValue $ws "D45" "0.5853"
Set-TextValue $ws "E45" "  -4.61%  "
Set-TextValue $ws "D46" "12.78"
Set-TextValue $ws "E46" "  -2.90%  "
Set-TextValue $ws "D47" "3.623"
Set-TextValue $ws "E47" "  -3.24%  "
Set-TextValue $ws "D48" "1.980"
Set-TextValue $ws "E48" "  -3.77%  "
Set-TextValue $ws "E49" "  -3.06%  "
Set-TextValue $ws "D50" "1.197"
Set-TextValue $ws "E50" "  -3.74%  "
Set-TextValue $ws "D51" "1.153"
Set-TextValue $ws "E51" "  +0.59%  "
